$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The investor KYC sheet previously tracked an imported "user" (first
# name, last name, email + a mailto hyperlink, phone, and a yes/no
# "send confirmation email" flag). That whole identity is being
# removed, leaving PAN/Address/Bank Account/IFSC Code plus a new
# simple "Verified" flag.

# Drop the mailto hyperlinks that lived on the Email column before the
# columns shift underneath them.
$ws.Hyperlinks.Delete()

# Remove "First Name", "Last Name", "Email" and "Phone" (C:F) and the
# now-unused trailing width-only columns (old H:J, which become blank
# once the tail columns collapse into C:G).
$ws.Range("C:F,H:J").EntireColumn.Delete()

# Replace the old "Send Confirmation Email" column (now G) with a new
# "Verified" column, every investor marked Yes.
$ws.Range("G1").Value = "Verified"
$ws.Range("G2").Value = "Yes"
$ws.Range("G3").Value = "Yes"

# Match the saved selection left behind by the edit.
$ws.Range("G4").Select()

# The Hyperlink cell style is no longer referenced by any cell; drop it.
$wb.Styles.Item("Hyperlink").Delete()
